$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "..._old" -> "..._FV2304", "..._new" -> "..._FV2310"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $newVal = $val
        if ($val.EndsWith("_old")) {
            $newVal = $val.Substring(0, $val.Length - 4) + "_FV2304"
        } elseif ($val.EndsWith("_new")) {
            $newVal = $val.Substring(0, $val.Length - 4) + "_FV2310"
        }
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# 2. Turn the data range into an Excel Table ("ListObject") with an AutoFilter
$rng = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (pane split after row 1)
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true

Write-Host "Edit complete"
